# Weekly update: a new price record was reported for "Vega Modelo de Temuco"
# / "Bruselas (repollito)" (Provincia de Quillota origin), dated 2022-08-15
# (serial 44775). It is inserted as the new row 44, pushing every existing
# record from row 44 onward down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44; all rows 44..105 shift to 45..106.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record's data.
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44775
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112035
$ws.Range("G44").Value = "Bruselas (repollito)"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 25
$ws.Range("K44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("M44").Value = 25000
$ws.Range("N44").Value = "`$/malla 10 kilos"
$ws.Range("O44").Value = "Provincia de Quillota"
$ws.Range("P44").Value = 2500
$ws.Range("Q44").Value = 10
$ws.Range("R44").Value = "Hortaliza"
